$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1. Row 34: responsible person changes from "Coen" to "Kelly?"
# ---------------------------------------------------------------------
$ws1.Range("J34").Value = "Kelly?"

# ---------------------------------------------------------------------
# 2. Row 46: add a note in column K (plain/default formatting).
# ---------------------------------------------------------------------
$ws1.Range("K46").Value = "Question Francesco"

# ---------------------------------------------------------------------
# 3. New bold header row (row 6): add Index/ISIN/Name/Type in C6:F6 and
#    bold the existing Responsible/Done?/DESC headers in J6:L6.
# ---------------------------------------------------------------------
$ws1.Range("C6").Value = "Index"
$ws1.Range("D6").Value = "ISIN"
$ws1.Range("E6").Value = "Name"
$ws1.Range("F6").Value = "Type"

$ws1.Range("C6:F6").Font.Bold = $true
$ws1.Range("C6:F6").Font.Name = "Verdana"
$ws1.Range("C6:F6").Font.Size = 10

$ws1.Range("J6:L6").Font.Bold = $true
$ws1.Range("J6:L6").Font.Name = "Verdana"
$ws1.Range("J6:L6").Font.Size = 10

# ---------------------------------------------------------------------
# 4. Add "Price" to column F for every data row (34-51), matching the
#    formatting already used by the other cells on those rows (style
#    copied from C34 so no stray new style gets created).
# ---------------------------------------------------------------------
$ws1.Range("C34").Copy()
$ws1.Range("F34:F51").PasteSpecial(-4122)
$ws1.Range("F34:F51").Value = "Price"

# ---------------------------------------------------------------------
# 5. Mark rows 35-37 as Done in column K (style copied from K29, an
#    existing "Done" cell with the same formatting).
# ---------------------------------------------------------------------
$ws1.Range("K29").Copy()
$ws1.Range("K35:K37").PasteSpecial(-4122)
$ws1.Range("K35:K37").Value = "Done"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6. Hide columns G, H and I (kept but no longer shown) and add a new
#    column K sized for the new "Done?" notes.
# ---------------------------------------------------------------------
$ws1.Columns("G").Hidden = $true
$ws1.Columns("H").Hidden = $true
$ws1.Columns("I").Hidden = $true
$ws1.Columns("K").ColumnWidth = 16.71

# ---------------------------------------------------------------------
# 7. Sheet2: record a selection on B40 without leaving it as the active
#    sheet (Sheet1 stays the active/tabbed sheet as before).
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B40").Select()
$ws1.Activate()
$ws1.Range("K46").Select()
